# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Tonberry_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1238.2174  # ALC!H15: 1198.5834 -> 1238.2174
$ws.Cells.Item(15, 9).Value = 1238.2174  # ALC!I15: 1198.5834 -> 1238.2174
$ws.Cells.Item(15, 11).Value = 3714.6522  # ALC!K15: 3595.7502 -> 3714.6522
$ws.Cells.Item(15, 13).Value = -3545.6522  # ALC!M15: -3426.7502 -> -3545.6522

$ws.Cells.Item(18, 8).Value = 16816.334  # ALC!H18: 16179.6 -> 16816.334

$ws.Cells.Item(19, 8).Value = 771367.44  # ALC!H19: 771681.3 -> 771367.4399999999
$ws.Cells.Item(19, 9).Value = 2000819.8  # ALC!I19: 2001295.8 -> 2000819.8
$ws.Cells.Item(19, 10).Value = 2959.75  # ALC!J19: 3172.25 -> 2959.75
$ws.Cells.Item(19, 11).Value = 2000819.8  # ALC!K19: 2001295.8 -> 2000819.8
$ws.Cells.Item(19, 12).Value = 2959.75  # ALC!L19: 3172.25 -> 2959.75
$ws.Cells.Item(19, 13).Value = -2000644.8  # ALC!M19: -2001120.8 -> -2000644.8
$ws.Cells.Item(19, 14).Value = -3309.75  # ALC!N19: -3522.25 -> -3309.75

$ws.Cells.Item(132, 8).Value = 1197.4166  # ALC!H132: 1224.4546 -> 1197.4166
$ws.Cells.Item(132, 9).Value = 1169.909  # ALC!I132: 1196.9 -> 1169.909
$ws.Cells.Item(132, 11).Value = 3509.727  # ALC!K132: 3590.7 -> 3509.727
$ws.Cells.Item(132, 13).Value = -979.7270000000003  # ALC!M132: -1060.7 -> -979.7270000000003

$ws.Cells.Item(137, 8).Value = 2432.4443  # ALC!H137: 2366.2222 -> 2432.4443
$ws.Cells.Item(137, 9).Value = 1778.6  # ALC!I137: 1882.8334 -> 1778.6
$ws.Cells.Item(137, 10).Value = 3249.75  # ALC!J137: 3333 -> 3249.75
$ws.Cells.Item(137, 11).Value = 5335.799999999999  # ALC!K137: 5648.5002 -> 5335.799999999999
$ws.Cells.Item(137, 12).Value = 9749.25  # ALC!L137: 9999 -> 9749.25
$ws.Cells.Item(137, 13).Value = -2785.799999999999  # ALC!M137: -3098.5002 -> -2785.799999999999
$ws.Cells.Item(137, 14).Value = -14849.25  # ALC!N137: -15099 -> -14849.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 265395.75  # ARM!H2: 265428.9 -> 265395.75
$ws.Cells.Item(2, 9).Value = 347894.62  # ARM!I2: 347938 -> 347894.62
$ws.Cells.Item(2, 10).Value = 1399.4  # ARM!J2: 1399.8 -> 1399.4
$ws.Cells.Item(2, 11).Value = 347894.62  # ARM!K2: 347938 -> 347894.62
$ws.Cells.Item(2, 12).Value = 1399.4  # ARM!L2: 1399.8 -> 1399.4
$ws.Cells.Item(2, 13).Value = -347781.62  # ARM!M2: -347825 -> -347781.62
$ws.Cells.Item(2, 14).Value = -1625.4  # ARM!N2: -1625.8 -> -1625.4

$ws.Cells.Item(32, 8).Value = 3308.7  # ARM!H32: 2781.6 -> 3308.7
$ws.Cells.Item(32, 9).Value = 3308.7  # ARM!I32: 2781.6 -> 3308.7
$ws.Cells.Item(32, 11).Value = 3308.7  # ARM!K32: 2781.6 -> 3308.7
$ws.Cells.Item(32, 13).Value = -3021.7  # ARM!M32: -2494.6 -> -3021.7

$ws.Cells.Item(45, 8).Value = 6001307  # ARM!H45: 6429921.5 -> 6001307
$ws.Cells.Item(45, 10).Value = 1611.6666  # ARM!J45: 1725.625 -> 1611.6666
$ws.Cells.Item(45, 12).Value = 1611.6666  # ARM!L45: 1725.625 -> 1611.6666
$ws.Cells.Item(45, 14).Value = -2365.6666  # ARM!N45: -2479.625 -> -2365.6666

$ws.Cells.Item(61, 8).Value = 22929.36  # ARM!H61: 22929.64 -> 22929.36
$ws.Cells.Item(61, 9).Value = 35100.168  # ARM!I61: 36583.26 -> 35100.168
$ws.Cells.Item(61, 10).Value = 3456.0667  # ARM!J61: 3302.5625 -> 3456.0667
$ws.Cells.Item(61, 11).Value = 35100.168  # ARM!K61: 36583.26 -> 35100.168
$ws.Cells.Item(61, 12).Value = 3456.0667  # ARM!L61: 3302.5625 -> 3456.0667
$ws.Cells.Item(61, 13).Value = -34888.168  # ARM!M61: -36371.26 -> -34888.168
$ws.Cells.Item(61, 14).Value = -3880.0667  # ARM!N61: -3726.5625 -> -3880.0667

$ws.Cells.Item(74, 8).Value = 912.1667  # ARM!H74: 962.1667 -> 912.1667
$ws.Cells.Item(74, 9).Value = 783.6111  # ARM!I74: 794.5714 -> 783.6111
$ws.Cells.Item(74, 10).Value = 1683.5  # ARM!J74: 1800.1428 -> 1683.5
$ws.Cells.Item(74, 11).Value = 783.6111  # ARM!K74: 794.5714 -> 783.6111
$ws.Cells.Item(74, 12).Value = 1683.5  # ARM!L74: 1800.1428 -> 1683.5
$ws.Cells.Item(74, 13).Value = 90.38890000000004  # ARM!M74: 79.42859999999996 -> 90.38890000000004
$ws.Cells.Item(74, 14).Value = -3431.5  # ARM!N74: -3548.1428 -> -3431.5

$ws.Cells.Item(77, 8).Value = 912.1667  # ARM!H77: 962.1667 -> 912.1667
$ws.Cells.Item(77, 9).Value = 783.6111  # ARM!I77: 794.5714 -> 783.6111
$ws.Cells.Item(77, 10).Value = 1683.5  # ARM!J77: 1800.1428 -> 1683.5
$ws.Cells.Item(77, 11).Value = 3918.0555  # ARM!K77: 3972.857 -> 3918.0555
$ws.Cells.Item(77, 12).Value = 8417.5  # ARM!L77: 9000.714 -> 8417.5
$ws.Cells.Item(77, 13).Value = 449.9445000000001  # ARM!M77: 395.143 -> 449.9445000000001
$ws.Cells.Item(77, 14).Value = -17153.5  # ARM!N77: -17736.714 -> -17153.5

$ws.Cells.Item(116, 8).Value = 265395.75  # ARM!H116: 265428.9 -> 265395.75
$ws.Cells.Item(116, 9).Value = 347894.62  # ARM!I116: 347938 -> 347894.62
$ws.Cells.Item(116, 10).Value = 1399.4  # ARM!J116: 1399.8 -> 1399.4
$ws.Cells.Item(116, 11).Value = 347894.62  # ARM!K116: 347938 -> 347894.62
$ws.Cells.Item(116, 12).Value = 1399.4  # ARM!L116: 1399.8 -> 1399.4
$ws.Cells.Item(116, 13).Value = -345600.62  # ARM!M116: -345644 -> -345600.62
$ws.Cells.Item(116, 14).Value = -5987.4  # ARM!N116: -5987.8 -> -5987.4

$ws.Cells.Item(122, 8).Value = 1509.5217  # ARM!H122: 1317.5518 -> 1509.5217
$ws.Cells.Item(122, 9).Value = 1524.9375  # ARM!I122: 1344.5 -> 1524.9375
$ws.Cells.Item(122, 10).Value = 1474.2858  # ARM!J122: 1257.6666 -> 1474.2858
$ws.Cells.Item(122, 11).Value = 4574.8125  # ARM!K122: 4033.5 -> 4574.8125
$ws.Cells.Item(122, 12).Value = 4422.857400000001  # ARM!L122: 3772.9998 -> 4422.857400000001
$ws.Cells.Item(122, 13).Value = -2124.8125  # ARM!M122: -1583.5 -> -2124.8125
$ws.Cells.Item(122, 14).Value = -9322.8574  # ARM!N122: -8672.9998 -> -9322.857400000001

$ws.Cells.Item(132, 8).Value = 2262.0908  # ARM!H132: 2332.2144 -> 2262.0908
$ws.Cells.Item(132, 9).Value = 2024.45  # ARM!I132: 2161.6667 -> 2024.45
$ws.Cells.Item(132, 11).Value = 6073.35  # ARM!K132: 6485.000100000001 -> 6073.35
$ws.Cells.Item(132, 13).Value = -3543.35  # ARM!M132: -3955.000100000001 -> -3543.35

$ws.Cells.Item(136, 8).Value = 22929.36  # ARM!H136: 22929.64 -> 22929.36
$ws.Cells.Item(136, 9).Value = 35100.168  # ARM!I136: 36583.26 -> 35100.168
$ws.Cells.Item(136, 10).Value = 3456.0667  # ARM!J136: 3302.5625 -> 3456.0667
$ws.Cells.Item(136, 11).Value = 105300.504  # ARM!K136: 109749.78 -> 105300.504
$ws.Cells.Item(136, 12).Value = 10368.2001  # ARM!L136: 9907.6875 -> 10368.2001
$ws.Cells.Item(136, 13).Value = -102750.504  # ARM!M136: -107199.78 -> -102750.504
$ws.Cells.Item(136, 14).Value = -15468.2001  # ARM!N136: -15007.6875 -> -15468.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 265395.75  # BSM!H3: 265428.9 -> 265395.75
$ws.Cells.Item(3, 9).Value = 347894.62  # BSM!I3: 347938 -> 347894.62
$ws.Cells.Item(3, 10).Value = 1399.4  # BSM!J3: 1399.8 -> 1399.4
$ws.Cells.Item(3, 11).Value = 347894.62  # BSM!K3: 347938 -> 347894.62
$ws.Cells.Item(3, 12).Value = 1399.4  # BSM!L3: 1399.8 -> 1399.4
$ws.Cells.Item(3, 13).Value = -347780.62  # BSM!M3: -347824 -> -347780.62
$ws.Cells.Item(3, 14).Value = -1627.4  # BSM!N3: -1627.8 -> -1627.4

$ws.Cells.Item(20, 8).Value = 3609  # BSM!H20: 3066.5715 -> 3609
$ws.Cells.Item(20, 9).Value = 1756.3334  # BSM!I20: 1623 -> 1756.3334
$ws.Cells.Item(20, 10).Value = 4998.5  # BSM!J20: 4149.25 -> 4998.5
$ws.Cells.Item(20, 11).Value = 1756.3334  # BSM!K20: 1623 -> 1756.3334
$ws.Cells.Item(20, 12).Value = 4998.5  # BSM!L20: 4149.25 -> 4998.5
$ws.Cells.Item(20, 13).Value = -1509.3334  # BSM!M20: -1376 -> -1509.3334
$ws.Cells.Item(20, 14).Value = -5492.5  # BSM!N20: -4643.25 -> -5492.5

$ws.Cells.Item(86, 8).Value = 800979.8  # BSM!H86: 667929 -> 800979.8
$ws.Cells.Item(86, 9).Value = 1449.5  # BSM!I86: 1858 -> 1449.5
$ws.Cells.Item(86, 11).Value = 1449.5  # BSM!K86: 1858 -> 1449.5
$ws.Cells.Item(86, 13).Value = -326.5  # BSM!M86: -735 -> -326.5

$ws.Cells.Item(89, 8).Value = 800979.8  # BSM!H89: 667929 -> 800979.8
$ws.Cells.Item(89, 9).Value = 1449.5  # BSM!I89: 1858 -> 1449.5
$ws.Cells.Item(89, 11).Value = 7247.5  # BSM!K89: 9290 -> 7247.5
$ws.Cells.Item(89, 13).Value = -1631.5  # BSM!M89: -3674 -> -1631.5

$ws.Cells.Item(134, 8).Value = 3579  # BSM!H134: 3663.348 -> 3579
$ws.Cells.Item(134, 9).Value = 2879.875  # BSM!I134: 3014.5334 -> 2879.875
$ws.Cells.Item(134, 10).Value = 5177  # BSM!J134: 4879.875 -> 5177
$ws.Cells.Item(134, 11).Value = 8639.625  # BSM!K134: 9043.600199999999 -> 8639.625
$ws.Cells.Item(134, 12).Value = 15531  # BSM!L134: 14639.625 -> 15531
$ws.Cells.Item(134, 13).Value = -6104.625  # BSM!M134: -6508.600199999999 -> -6104.625
$ws.Cells.Item(134, 14).Value = -20601  # BSM!N134: -19709.625 -> -20601

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2315.0588  # CRP!H31: 2303.0833 -> 2315.0588
$ws.Cells.Item(31, 10).Value = 2501.8235  # CRP!J31: 2459.4736 -> 2501.8235
$ws.Cells.Item(31, 12).Value = 2501.8235  # CRP!L31: 2459.4736 -> 2501.8235
$ws.Cells.Item(31, 14).Value = -3091.8235  # CRP!N31: -3049.4736 -> -3091.8235

$ws.Cells.Item(34, 8).Value = 2315.0588  # CRP!H34: 2303.0833 -> 2315.0588
$ws.Cells.Item(34, 10).Value = 2501.8235  # CRP!J34: 2459.4736 -> 2501.8235
$ws.Cells.Item(34, 12).Value = 2501.8235  # CRP!L34: 2459.4736 -> 2501.8235
$ws.Cells.Item(34, 14).Value = -2905.8235  # CRP!N34: -2863.4736 -> -2905.8235

$ws.Cells.Item(105, 8).Value = 777.25  # CRP!H105: 821.1429000000001 -> 777.25
$ws.Cells.Item(105, 9).Value = 789.5  # CRP!I105: 853.4 -> 789.5
$ws.Cells.Item(105, 11).Value = 789.5  # CRP!K105: 853.4 -> 789.5
$ws.Cells.Item(105, 13).Value = 957.5  # CRP!M105: 893.6 -> 957.5

$ws.Cells.Item(132, 8).Value = 1678.8889  # CRP!H132: 1719.4117 -> 1678.8889
$ws.Cells.Item(132, 9).Value = 1261.4  # CRP!I132: 1280.7858 -> 1261.4
$ws.Cells.Item(132, 11).Value = 3784.2  # CRP!K132: 3842.3574 -> 3784.2
$ws.Cells.Item(132, 13).Value = -1254.2  # CRP!M132: -1312.3574 -> -1254.2

$ws.Cells.Item(134, 8).Value = 1332.2812  # CRP!H134: 1332.3429 -> 1332.2812
$ws.Cells.Item(134, 9).Value = 1201.2222  # CRP!I134: 1194 -> 1201.2222
$ws.Cells.Item(134, 10).Value = 2040  # CRP!J134: 1885.7142 -> 2040
$ws.Cells.Item(134, 11).Value = 3603.6666  # CRP!K134: 3582 -> 3603.6666
$ws.Cells.Item(134, 12).Value = 6120  # CRP!L134: 5657.142599999999 -> 6120
$ws.Cells.Item(134, 13).Value = -1068.6666  # CRP!M134: -1047 -> -1068.6666
$ws.Cells.Item(134, 14).Value = -11190  # CRP!N134: -10727.1426 -> -11190

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 670.4138  # CUL!H5: 704.8077 -> 670.4138
$ws.Cells.Item(5, 9).Value = 593.95  # CUL!I5: 624.55554 -> 593.95
$ws.Cells.Item(5, 10).Value = 840.3333  # CUL!J5: 885.375 -> 840.3333
$ws.Cells.Item(5, 11).Value = 1781.85  # CUL!K5: 1873.66662 -> 1781.85
$ws.Cells.Item(5, 12).Value = 2520.9999  # CUL!L5: 2656.125 -> 2520.9999
$ws.Cells.Item(5, 13).Value = -1669.85  # CUL!M5: -1761.66662 -> -1669.85
$ws.Cells.Item(5, 14).Value = -2744.9999  # CUL!N5: -2880.125 -> -2744.9999

$ws.Cells.Item(117, 8).Value = 764.4286  # CUL!H117: 1145.2858 -> 764.4286
$ws.Cells.Item(117, 10).Value = 1651.5  # CUL!J117: 2984.5 -> 1651.5
$ws.Cells.Item(117, 12).Value = 4954.5  # CUL!L117: 8953.5 -> 4954.5
$ws.Cells.Item(117, 14).Value = -11838.5  # CUL!N117: -15837.5 -> -11838.5

$ws.Cells.Item(135, 8).Value = 670.4138  # CUL!H135: 704.8077 -> 670.4138
$ws.Cells.Item(135, 9).Value = 593.95  # CUL!I135: 624.55554 -> 593.95
$ws.Cells.Item(135, 10).Value = 840.3333  # CUL!J135: 885.375 -> 840.3333
$ws.Cells.Item(135, 11).Value = 5345.55  # CUL!K135: 5620.99986 -> 5345.55
$ws.Cells.Item(135, 12).Value = 7562.9997  # CUL!L135: 7968.375 -> 7562.9997
$ws.Cells.Item(135, 13).Value = -2810.55  # CUL!M135: -3085.99986 -> -2810.55
$ws.Cells.Item(135, 14).Value = -12632.9997  # CUL!N135: -13038.375 -> -12632.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 20000  # GSM!H34: 0 -> 20000
$ws.Cells.Item(34, 10).Value = 20000  # GSM!J34: 0 -> 20000
$ws.Cells.Item(34, 12).Value = 20000  # GSM!L34: 0 -> 20000
$ws.Cells.Item(34, 14).Value = -20536  # GSM!N34: None -> -20536

$ws.Cells.Item(70, 8).Value = 4176.1665  # GSM!H70: 4036.5715 -> 4176.1665
$ws.Cells.Item(70, 9).Value = 3849.75  # GSM!I70: 3749.75 -> 3849.75
$ws.Cells.Item(70, 10).Value = 4339.375  # GSM!J70: 4151.3 -> 4339.375
$ws.Cells.Item(70, 11).Value = 3849.75  # GSM!K70: 3749.75 -> 3849.75
$ws.Cells.Item(70, 12).Value = 4339.375  # GSM!L70: 4151.3 -> 4339.375
$ws.Cells.Item(70, 13).Value = -3579.75  # GSM!M70: -3479.75 -> -3579.75
$ws.Cells.Item(70, 14).Value = -4879.375  # GSM!N70: -4691.3 -> -4879.375

$ws.Cells.Item(73, 8).Value = 4176.1665  # GSM!H73: 4036.5715 -> 4176.1665
$ws.Cells.Item(73, 9).Value = 3849.75  # GSM!I73: 3749.75 -> 3849.75
$ws.Cells.Item(73, 10).Value = 4339.375  # GSM!J73: 4151.3 -> 4339.375
$ws.Cells.Item(73, 11).Value = 3849.75  # GSM!K73: 3749.75 -> 3849.75
$ws.Cells.Item(73, 12).Value = 4339.375  # GSM!L73: 4151.3 -> 4339.375
$ws.Cells.Item(73, 13).Value = -2913.75  # GSM!M73: -2813.75 -> -2913.75
$ws.Cells.Item(73, 14).Value = -6211.375  # GSM!N73: -6023.3 -> -6211.375

$ws.Cells.Item(76, 8).Value = 20000  # GSM!H76: 0 -> 20000
$ws.Cells.Item(76, 10).Value = 20000  # GSM!J76: 0 -> 20000
$ws.Cells.Item(76, 12).Value = 20000  # GSM!L76: 0 -> 20000
$ws.Cells.Item(76, 14).Value = -20630  # GSM!N76: None -> -20630

$ws.Cells.Item(79, 8).Value = 20000  # GSM!H79: 0 -> 20000
$ws.Cells.Item(79, 10).Value = 20000  # GSM!J79: 0 -> 20000
$ws.Cells.Item(79, 12).Value = 20000  # GSM!L79: 0 -> 20000
$ws.Cells.Item(79, 14).Value = -22184  # GSM!N79: None -> -22184

$ws.Cells.Item(122, 8).Value = 1476.5  # GSM!H122: 1640 -> 1476.5
$ws.Cells.Item(122, 9).Value = 1340  # GSM!I122: 1425 -> 1340
$ws.Cells.Item(122, 10).Value = 1704  # GSM!J122: 2500 -> 1704
$ws.Cells.Item(122, 11).Value = 4020  # GSM!K122: 4275 -> 4020
$ws.Cells.Item(122, 12).Value = 5112  # GSM!L122: 7500 -> 5112
$ws.Cells.Item(122, 13).Value = -1570  # GSM!M122: -1825 -> -1570
$ws.Cells.Item(122, 14).Value = -10012  # GSM!N122: -12400 -> -10012

$ws.Cells.Item(132, 8).Value = 1104007.1  # GSM!H132: 1136454.4 -> 1104007.1
$ws.Cells.Item(132, 9).Value = 1484866.5  # GSM!I132: 1544229.1 -> 1484866.5
$ws.Cells.Item(132, 11).Value = 4454599.5  # GSM!K132: 4632687.300000001 -> 4454599.5
$ws.Cells.Item(132, 13).Value = -4452069.5  # GSM!M132: -4630157.300000001 -> -4452069.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 20000520  # LTW!H55: 25000538 -> 20000520
$ws.Cells.Item(55, 9).Value = 33333684  # LTW!I55: 50000300 -> 33333684
$ws.Cells.Item(55, 11).Value = 33333684  # LTW!K55: 50000300 -> 33333684
$ws.Cells.Item(55, 13).Value = -33333511  # LTW!M55: -50000127 -> -33333511

$ws.Cells.Item(61, 8).Value = 2333.5293  # LTW!H61: 2198.3684 -> 2333.5293
$ws.Cells.Item(61, 9).Value = 2244.2307  # LTW!I61: 2084.9333 -> 2244.2307
$ws.Cells.Item(61, 11).Value = 2244.2307  # LTW!K61: 2084.9333 -> 2244.2307
$ws.Cells.Item(61, 13).Value = -2042.2307  # LTW!M61: -1882.9333 -> -2042.2307

$ws.Cells.Item(113, 8).Value = 2333.5293  # LTW!H113: 2198.3684 -> 2333.5293
$ws.Cells.Item(113, 9).Value = 2244.2307  # LTW!I113: 2084.9333 -> 2244.2307
$ws.Cells.Item(113, 11).Value = 2244.2307  # LTW!K113: 2084.9333 -> 2244.2307
$ws.Cells.Item(113, 13).Value = -74.23070000000007  # LTW!M113: 85.06669999999986 -> -74.23070000000007

$ws.Cells.Item(132, 8).Value = 4946.4  # LTW!H132: 5277.278 -> 4946.4
$ws.Cells.Item(132, 9).Value = 3283.182  # LTW!I132: 3575.3333 -> 3283.182
$ws.Cells.Item(132, 11).Value = 9849.545999999998  # LTW!K132: 10725.9999 -> 9849.545999999998
$ws.Cells.Item(132, 13).Value = -7319.545999999998  # LTW!M132: -8195.999899999999 -> -7319.545999999998

$ws.Cells.Item(136, 8).Value = 1146  # LTW!H136: 1091.6 -> 1146
$ws.Cells.Item(136, 9).Value = 1146  # LTW!I136: 1091.6 -> 1146
$ws.Cells.Item(136, 11).Value = 3438  # LTW!K136: 3274.8 -> 3438
$ws.Cells.Item(136, 13).Value = -888  # LTW!M136: -724.7999999999997 -> -888

$ws.Cells.Item(137, 8).Value = 55164.5  # LTW!H137: 0 -> 55164.5
$ws.Cells.Item(137, 10).Value = 55164.5  # LTW!J137: 0 -> 55164.5
$ws.Cells.Item(137, 12).Value = 55164.5  # LTW!L137: 0 -> 55164.5
$ws.Cells.Item(137, 14).Value = -65364.5  # LTW!N137: None -> -65364.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 0  # WVR!H4: 1074 -> 0
$ws.Cells.Item(4, 10).Value = 0  # WVR!J4: 1074 -> 0
$ws.Cells.Item(4, 12).Value = 0  # WVR!L4: 1074 -> 0
$ws.Cells.Item(4, 14).Value = ""  # WVR!N4: -1300 -> (cleared)

$ws.Cells.Item(5, 8).Value = 9900  # WVR!H5: 3179.8 -> 9900
$ws.Cells.Item(5, 9).Value = 0  # WVR!I5: 1199 -> 0
$ws.Cells.Item(5, 10).Value = 9900  # WVR!J5: 3675 -> 9900
$ws.Cells.Item(5, 11).Value = 0  # WVR!K5: 1199 -> 0
$ws.Cells.Item(5, 12).Value = 9900  # WVR!L5: 3675 -> 9900
$ws.Cells.Item(5, 13).Value = ""  # WVR!M5: -1087 -> (cleared)
$ws.Cells.Item(5, 14).Value = -10124  # WVR!N5: -3899 -> -10124

$ws.Cells.Item(122, 8).Value = 43889.844  # WVR!H122: 36412.477 -> 43889.844
$ws.Cells.Item(122, 9).Value = 58486.215  # WVR!I122: 48323.06 -> 58486.215
$ws.Cells.Item(122, 10).Value = 3020  # WVR!J122: 2665.8333 -> 3020
$ws.Cells.Item(122, 11).Value = 175458.645  # WVR!K122: 144969.18 -> 175458.645
$ws.Cells.Item(122, 12).Value = 9060  # WVR!L122: 7997.499899999999 -> 9060
$ws.Cells.Item(122, 13).Value = -173008.645  # WVR!M122: -142519.18 -> -173008.645
$ws.Cells.Item(122, 14).Value = -13960  # WVR!N122: -12897.4999 -> -13960

$ws.Cells.Item(126, 8).Value = 7324.0527  # WVR!H126: 7324.579 -> 7324.0527
$ws.Cells.Item(126, 9).Value = 8178.467  # WVR!I126: 8179.1333 -> 8178.467
$ws.Cells.Item(126, 11).Value = 24535.401  # WVR!K126: 24537.3999 -> 24535.401
$ws.Cells.Item(126, 13).Value = -22065.401  # WVR!M126: -22067.3999 -> -22065.401

$ws.Cells.Item(132, 8).Value = 1899.814  # WVR!H132: 1921.238 -> 1899.814
$ws.Cells.Item(132, 9).Value = 1745.8125  # WVR!I132: 1769.871 -> 1745.8125
$ws.Cells.Item(132, 11).Value = 5237.4375  # WVR!K132: 5309.613 -> 5237.4375
$ws.Cells.Item(132, 13).Value = -2707.4375  # WVR!M132: -2779.613 -> -2707.4375

$ws.Cells.Item(136, 8).Value = 15433532  # WVR!H136: 15874476 -> 15433532
$ws.Cells.Item(136, 9).Value = 21368454  # WVR!I136: 22223174 -> 21368454
$ws.Cells.Item(136, 11).Value = 64105362  # WVR!K136: 66669522 -> 64105362
$ws.Cells.Item(136, 13).Value = -64102812  # WVR!M136: -66666972 -> -64102812
